$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("steel_prim")
$ws.Range("B13").Value = 79.2039275472

$ws = $wb.Worksheets.Item("steel_sec")
$ws.Range("B2").Value = 3137.3333638953
$ws.Range("B3").Value = 662.8772678271
$ws.Range("B4").Value = 187.1929792754
$ws.Range("B6").Value = 14745.4804747338
$ws.Range("B8").Value = 1686.9328934361
$ws.Range("B9").Value = 12000.1768022534
$ws.Range("B10").Value = 10578.4163151474
$ws.Range("B11").Value = 74.52604210760001
$ws.Range("B12").Value = 19166.1751082983
$ws.Range("B14").Value = 2918.2085010061
$ws.Range("B15").Value = 351.581624267
$ws.Range("B17").Value = 780.4161820822
$ws.Range("B18").Value = 6304.7448647049
$ws.Range("B19").Value = 2229.0851216163
$ws.Range("B20").Value = 1249.3365025614
$ws.Range("B21").Value = 672.0149639492
$ws.Range("B22").Value = 401.5685422886
$ws.Range("B23").Value = 1550.7885438499
$ws.Range("B24").Value = 1617.5105020949
$ws.Range("B25").Value = 1491.5147196376
$ws.Range("B26").Value = 681.3453594332
$ws.Range("B27").Value = 1487.1731822742
$ws.Range("B28").Value = 135.1748448513
$ws.Range("B29").Value = 276.6780234494
$ws.Range("B31").Value = 203

$ws = $wb.Worksheets.Item("alu_prim")
$ws.Range("B6").Value = 484.6363054935
$ws.Range("B8").Value = 169.5033804624
$ws.Range("B9").Value = 328.5140390669
$ws.Range("B10").Value = 719.4540185677
$ws.Range("B16").Value = 83.38897798230001
$ws.Range("B20").Value = 286.0628334636
$ws.Range("B21").Value = 71.7067978
$ws.Range("B22").Value = 209.4584584959
$ws.Range("B24").Value = 124.6379539807
$ws.Range("B25").Value = 123.4637289649
$ws.Range("B26").Value = 1482.6349270892
$ws.Range("B28").Value = 45.0420447981
$ws.Range("B32").Value = 133.7
$ws.Range("B33").Value = 985.3530527276

$ws = $wb.Worksheets.Item("chlorine")
$ws.Range("B2").Value = 915.4614793683
$ws.Range("B4").Value = 121.6660190167
$ws.Range("B6").Value = 4446.2757975824
$ws.Range("B7").Value = 16.411172949
$ws.Range("B8").Value = 20.3124846355
$ws.Range("B9").Value = 433.657834659
$ws.Range("B10").Value = 2007.6908315921
$ws.Range("B12").Value = 141.4187606185
$ws.Range("B15").Value = 430.0080968944
$ws.Range("B16").Value = 725.3109006993
$ws.Range("B17").Value = 61.5569201487
$ws.Range("B18").Value = 400.3995477578
$ws.Range("B19").Value = 114.7067231633
$ws.Range("B20").Value = 155.4600549801
$ws.Range("B21").Value = 12.956235316
$ws.Range("B22").Value = 60.9537256277
$ws.Range("B23").Value = 78.7361214906
$ws.Range("B24").Value = 59.5760365719
$ws.Range("B25").Value = 449.5275538241
$ws.Range("B26").Value = 288.6122470802
$ws.Range("B27").Value = 16.0791044189

$ws = $wb.Worksheets.Item("paper")
$ws.Range("B2").Value = 2386.7350504722
$ws.Range("B3").Value = 350.8086148922
$ws.Range("B4").Value = 850.7451905665999
$ws.Range("B5").Value = 499.6364464216
$ws.Range("B6").Value = 25292.6666919066
$ws.Range("B7").Value = 134.8151374435
$ws.Range("B8").Value = 485.3765758082
$ws.Range("B9").Value = 6979.8852341307
$ws.Range("B10").Value = 15683.5795022683
$ws.Range("B11").Value = 454.883635958
$ws.Range("B12").Value = 10150.3121293548
$ws.Range("B13").Value = 40.0599767601
$ws.Range("B14").Value = 7.102246231
$ws.Range("B15").Value = 851.590154436
$ws.Range("B16").Value = 2862.1863440867
$ws.Range("B17").Value = 5562.0485442637
$ws.Range("B18").Value = 6505.8229817014
$ws.Range("B19").Value = 2198.1129850901
$ws.Range("B20").Value = 379.3635365228
$ws.Range("B21").Value = 818.9736998766
$ws.Range("B22").Value = 955.2138812278999
$ws.Range("B23").Value = 11801.1189834848
$ws.Range("B24").Value = 11927.0047822509
$ws.Range("B25").Value = 3873.4286425043
$ws.Range("B26").Value = 1227.9983719078
$ws.Range("B27").Value = 1432.6666783121
$ws.Range("B29").Value = 26.4365125137
$ws.Range("B31").Value = 493.3390385661
$ws.Range("B33").Value = 3.4172842242
$ws.Range("B34").Value = 290.3223524823
$ws.Range("B35").Value = 80.0779070815

$ws = $wb.Worksheets.Item("cement")
$ws.Range("B2").Value = 6531.1982775692
$ws.Range("B3").Value = 2126.274101738
$ws.Range("B4").Value = 3350.8914350681
$ws.Range("B5").Value = 1849.2452520753
$ws.Range("B6").Value = 30778.3853798394
$ws.Range("B7").Value = 6368.9859802827
$ws.Range("B8").Value = 8407.285384884301
$ws.Range("B9").Value = 23759.61301459
$ws.Range("B10").Value = 28990.4430808753
$ws.Range("B11").Value = 2842.6595158841
$ws.Range("B12").Value = 26657.0622826651
$ws.Range("B13").Value = 957.5806911955
$ws.Range("B14").Value = 1167.1168804713
$ws.Range("B15").Value = 2986.6320839199
$ws.Range("B16").Value = 2036.5266270393
$ws.Range("B17").Value = 4528.4900772175
$ws.Range("B18").Value = 22524.1970888768
$ws.Range("B19").Value = 6201.9623762503
$ws.Range("B20").Value = 7336.0704149343
$ws.Range("B21").Value = 973.402211502
$ws.Range("B22").Value = 3581.1399966184
$ws.Range("B23").Value = 1392.9889689929
$ws.Range("B24").Value = 2861.9945526084
$ws.Range("B25").Value = 9443.3921941013
$ws.Range("B26").Value = 1862.023932759
$ws.Range("B27").Value = 4395.1816139407
$ws.Range("B29").Value = 912.1466581925999
$ws.Range("B30").Value = 2357.3164799636
$ws.Range("B31").Value = 1817.9992497272
$ws.Range("B32").Value = 849.67981936
$ws.Range("B33").Value = 177.5086472811
$ws.Range("B34").Value = 1786.4214860609
$ws.Range("B35").Value = 528.7463960865
